$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7477669715881348
$ws.Range("B1").Value = 1.151355266571045
$ws.Range("C1").Value = 2.636869192123413
$ws.Range("D1").Value = 3.474196672439575
$ws.Range("E1").Value = 1.743422746658325
